$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "50.983.59"
$ws.Range("E2").Value = "  -2.07%  "
$ws.Range("D3").Value = "2.757.18"
$ws.Range("E3").Value = "  -1.11%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "352.48"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.77%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "107.59"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.65%  "
$ws.Range("E7").Value = "  -2.51%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.00"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.582"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.68%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "39.35"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.82%  "
$ws.Range("E11").Value = "  +3.68%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0833"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.53%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "19.78"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.73%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.50"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.71%  "
$ws.Range("D15").Value = "3.186.45"
$ws.Range("E15").Value = "  -1.17%  "
$ws.Range("D16").Value = "2.756.59"
$ws.Range("E16").Value = "  -1.69%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.929"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.06%  "
$ws.Range("D18").Value = "50.941.56"
$ws.Range("E18").Value = "  -1.90%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.63"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.55%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "3.06"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.45%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.02"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.48%  "
$ws.Range("D22").Value = "0.0₃0959"
$ws.Range("E22").Value = "  -2.05%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "69.50"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "264.67"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -3.19%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.71"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.23%  "
$ws.Range("E26").Value = "  +0.05%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "25.89"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.88%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.161"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +12.47%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "10.09"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.57%  "
$ws.Range("E30").Value = "  +0.63%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "51.66"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.67%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.06"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +5.52%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "34.25"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.57%  "
$ws.Range("E34").Value = "  -5.14%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.38"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.79%  "
$ws.Range("E36").Value = "  -1.45%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.00"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.02%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "18.25"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.30%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.13"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.61%  "
$ws.Range("E40").Value = "  -2.84%  "
$ws.Range("E41").Value = "  -0.83%  "
$ws.Range("E42").Value = "  -2.87%  "
$ws.Range("B43").Value = "Monero"
$ws.Range("C43").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "120.39"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.45%  "
$ws.Range("B44").Value = "EnergySwap"
$ws.Range("C44").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "21.98"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.28%  "
$ws.Range("E45").Value = "  -2.73%  "
$ws.Range("D46").Value = "2.083.36"
$ws.Range("E46").Value = "  +0.78%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.23"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.83%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.26"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.38%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.918"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.33%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "5.45"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -4.65%  "
$ws.Range("E51").Value = "  +3.93%  "
